$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 1407
$ws.Range("F9").Value = 10
$ws.Range("F11").Value = 645
$ws.Range("F13").Value = 117
$ws.Range("F14").Value = 1310
$ws.Range("F15").Value = 481
$ws.Range("F16").Value = 472
$ws.Range("F27").Value = 108
$ws.Range("F31").Value = 87

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 728
$ws.Range("F5").Value = 618
$ws.Range("F6").Value = 618
$ws.Range("F15").Value = 389
$ws.Range("F16").Value = 389
$ws.Range("F27").Value = 224

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 914
$ws.Range("F9").Value = 1140
$ws.Range("F10").Value = 264

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 914
$ws.Range("F10").Value = 1140
$ws.Range("F11").Value = 264
$ws.Range("F16").Value = 1407
$ws.Range("F17").Value = 618
$ws.Range("C19").Value = "上海·Coser新春拜年盛典"
$ws.Range("D19").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E19").Value = "2024.02.14 10:00-02.15 17:00"
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81588"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202402/FZxjoQH41706757037933.jpeg"
$ws.Range("C20").Value = "上海·原X铁X崩only"
$ws.Range("D20").Value = "澳门路168号 月星国际家居"
$ws.Range("E20").Value = "2024.02.14 10:30-02.14 16:30"
$ws.Range("F20").Value = 106
$ws.Range("G20").Value = 60
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81446"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202401/IIePRulM1706248855263.jpeg"
$ws.Range("C21").Value = "上海·奇卡波利动漫嘉年华"
$ws.Range("D21").Value = "申滨路36号 虹桥丽宝广场"
$ws.Range("E21").Value = "2024.02.14 10:00-02.14 17:00"
$ws.Range("F21").Value = 645
$ws.Range("G21").Value = 48.8
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=81260"
$ws.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202401/9OHovK2V1705978109130.jpeg"
$ws.Range("B22").Value = "2024-02-14"
$ws.Range("C22").Value = "上海·魔都COS漫展-情人节专场AM01"
$ws.Range("D22").Value = "澳门路168号月星家居(江宁路地铁站1号口步行420米) 月星广场"
$ws.Range("E22").Value = "2024.02.14 10:00-02.14 16:00"
$ws.Range("F22").Value = 132
$ws.Range("G22").Value = 49
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=80691"
$ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202401/aSdjV6Kw1704868345679.jpeg"
$ws.Range("B23").Value = "2024-02-15"
$ws.Range("C23").Value = "上海·飘起来吧魔法泡泡-魔术表演"
$ws.Range("D23").Value = "曹杨路1888号 上海露边社·演艺空间"
$ws.Range("E23").Value = "2024.02.15 19:00-03.03 20:10"
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 88
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=81524"
$ws.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202401/tls18D0J1706599640356.png"
$ws.Range("B24").Value = "2024-02-16"
$ws.Range("C24").Value = "上海·次元裂缝-X 新年anikura派对"
$ws.Range("D24").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E24").Value = "2024.02.16 14:00-02.16 19:00"
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=81314"
$ws.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202401/OrhHWKdR1706062360956.jpeg"
$ws.Range("B25").Value = "2024-02-17"
$ws.Range("C25").Value = "上海·少女番only2.0"
$ws.Range("D25").Value = "营口路699号(黄兴公园地铁站2号口旁) 花嫁丽舍"
$ws.Range("E25").Value = "2024.02.17 10:00-02.17 17:00"
$ws.Range("F25").Value = 481
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=81148"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202401/j6eEZ18S1705657346664.jpeg"
$ws.Range("F26").Value = 472
$ws.Range("F35").Value = 108
$ws.Range("F39").Value = 389
$ws.Range("F42").Value = 87
$ws.Range("F47").Value = 224
$ws.Range("B49").Value = "2024-05-01"
$ws.Range("C49").Value = "上海·S·CGE动漫游戏嘉年华"
$ws.Range("D49").Value = "军工路1076号 纪希片场(秀场)"
$ws.Range("E49").Value = "2024.05.01 10:00-05.02 17:00"
$ws.Range("F49").Value = 28
$ws.Range("G49").Value = 70
$ws.Range("H49").Value = "https://show.bilibili.com/platform/detail.html?id=81204"
$ws.Range("I49").Value = "//i0.hdslb.com/bfs/openplatform/202401/nbFRULYe1705904589212.jpeg"
$ws.Range("C50").Value = "上海·魔都野良神only"
$ws.Range("D50").Value = "南京东路830号 第一百货"
$ws.Range("E50").Value = "2024.05.01 10:00-05.01 17:00"
$ws.Range("F50").Value = 241
$ws.Range("G50").Value = 79
$ws.Range("H50").Value = "https://show.bilibili.com/platform/detail.html?id=80321"
$ws.Range("I50").Value = "//i2.hdslb.com/bfs/openplatform/202401/KBlb0enU1704358750268.jpeg"
$ws.Range("B51").Value = "2024-05-05"
$ws.Range("C51").Value = "上海·灌篮高手--青春永不散场"
$ws.Range("D51").Value = "逸仙路1328弄 新业坊"
$ws.Range("E51").Value = "2024.05.05 10:00-05.05 17:00"
$ws.Range("F51").Value = 25
$ws.Range("G51").Value = 89
$ws.Range("H51").Value = "https://show.bilibili.com/platform/detail.html?id=80835"
$ws.Range("I51").Value = "//i1.hdslb.com/bfs/openplatform/202401/hdaVclFC1705301931054.jpeg"
